# sekhar - changes made to salaries
#
# 1. Insert a new title row at the top of the sheet:
#       "Bank Details and Netpay of Employees"
#    merged across A1:F1, styled Arial 15 bold blue, centered.
# 2. Widen column A so the long title is readable.
# 3. Update Sekhar Beri's Netpay (row 3) from 8701.46 -> 8700.73.
# 4. Update Pattabhi RamaRao Galidevara's Netpay (last row) from 11000.0 -> 10000.0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new header/title row, pushing everything else down one row ---
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Bank Details and Netpay of Employees"

# Merge the title across A1:F1
$ws.Range("A1:F1").Merge()

# Font: Arial, 15pt, bold, blue (COM colors are 0x00BBGGRR -> blue = 16711680)
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 15
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Color = 16711680

# Center the title
$ws.Range("A1").HorizontalAlignment = -4108

# Widen column A to comfortably fit the new title text
$ws.Columns.Item(1).ColumnWidth = 57.29

# --- Update the Netpay figures (rows shifted down by 1 after the insert) ---
# Row 3 = Sekhar Beri
$ws.Range("C3").Value = 8700.73
# Row 5 = Pattabhi RamaRao Galidevara
$ws.Range("C5").Value = 10000.0

Write-Output "Applied bank statement salary changes"
